# Update column F ("dSF") values on Sheet1 to match repulled/recomputed data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -7
$ws.Range("F7").Value = -5
$ws.Range("F11").Value = -4
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = -4
